# Fruta / hortaliza, semanal
# Insert 3 new weekly report rows at the top of the data block (row 25),
# pushing the existing rows 25-57 down to 28-60.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows before the current row 25.
$ws.Rows("25:27").Insert()

# Remember the date format used by the "Fecha" column (D) so the new rows
# keep the same look as the rest of the table.
$dateFmt = $ws.Cells.Item(28, 4).NumberFormat

# Columns A, B, C, E, F, G, H, I, J, K are constant for every record in this
# sheet (same market / product), so copy them from the row right below.
$constCols = @(1,2,3,5,6,7,8,9,10,11)

function Set-ConstantColumns($row) {
    foreach ($col in $constCols) {
        $ws.Cells.Item($row, $col).Value = $ws.Cells.Item(28, $col).Value2
    }
}

Set-ConstantColumns 25
Set-ConstantColumns 26
Set-ConstantColumns 27

# Row 25
$ws.Cells.Item(25, 4).Value = 45028
$ws.Cells.Item(25, 4).NumberFormat = $dateFmt
$ws.Cells.Item(25, 12).Value = "Especial"
$ws.Cells.Item(25, 13).Value = 330
$ws.Cells.Item(25, 14).Value = 13500
$ws.Cells.Item(25, 15).Value = 13500
$ws.Cells.Item(25, 16).Value = 13500
$ws.Cells.Item(25, 17).Value = "$/caja 15 kilos empedrada"
$ws.Cells.Item(25, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(25, 19).Value = 900
$ws.Cells.Item(25, 20).Value = 15

# Row 26
$ws.Cells.Item(26, 4).Value = 45028
$ws.Cells.Item(26, 4).NumberFormat = $dateFmt
$ws.Cells.Item(26, 12).Value = "Primera"
$ws.Cells.Item(26, 13).Value = 300
$ws.Cells.Item(26, 14).Value = 10500
$ws.Cells.Item(26, 15).Value = 10500
$ws.Cells.Item(26, 16).Value = 10500
$ws.Cells.Item(26, 17).Value = "$/caja 15 kilos empedrada"
$ws.Cells.Item(26, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(26, 19).Value = 700
$ws.Cells.Item(26, 20).Value = 15

# Row 27
$ws.Cells.Item(27, 4).Value = 45028
$ws.Cells.Item(27, 4).NumberFormat = $dateFmt
$ws.Cells.Item(27, 12).Value = "Segunda"
$ws.Cells.Item(27, 13).Value = 280
$ws.Cells.Item(27, 14).Value = 7500
$ws.Cells.Item(27, 15).Value = 7500
$ws.Cells.Item(27, 16).Value = 7500
$ws.Cells.Item(27, 17).Value = "$/caja 15 kilos empedrada"
$ws.Cells.Item(27, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(27, 19).Value = 500
$ws.Cells.Item(27, 20).Value = 15

Write-Host "Applied weekly Membrillo update: inserted rows 25-27"
